$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text stays as text (avoid Excel auto-converting numeric-looking
# strings like "59.026.38" or "1.00" into numbers/dates) by switching the
# Price/Volume columns to a text format before writing the new values, then
# restoring the original (Normal/General) style afterwards.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '59.026.38'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '2.572.39'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '570.86'
$ws.Range('E5').Value = '  +2.01%  '
$ws.Range('D6').Value = '142.66'
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  -0.65%  '
$ws.Range('D9').Value = '2.577.05'
$ws.Range('E9').Value = '  -1.93%  '
$ws.Range('E10').Value = '  -1.93%  '
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('E12').Value = '  +12.19%  '
$ws.Range('E13').Value = '  +2.28%  '
$ws.Range('D14').Value = '3.026.32'
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').Value = '59.083.54'
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').Value = '22.31'
$ws.Range('E16').Value = '  +5.52%  '
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('D18').Value = '2.577.24'
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('E19').Value = '  +1.27%  '
$ws.Range('D20').Value = '335.98'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').Value = '10.26'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '64.57'
$ws.Range('E24').Value = '  -2.30%  '
$ws.Range('D25').Value = '0.462'
$ws.Range('E25').Value = '  +7.81%  '
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -2.47%  '
$ws.Range('D28').Value = '7.28'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').Value = '0.0₃0779'
$ws.Range('E29').Value = '  +1.53%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').Value = '1.68'
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('D32').Value = '160.16'
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('D36').Value = '1.17'
$ws.Range('E36').Value = '  +2.37%  '
$ws.Range('E37').Value = '  -3.73%  '
$ws.Range('D38').Value = '0.873'
$ws.Range('E38').Value = '  -4.36%  '
$ws.Range('D39').Value = '37.42'
$ws.Range('E39').Value = '  +0.53%  '
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('D41').Value = '295.19'
$ws.Range('E41').Value = '  +3.49%  '
$ws.Range('D42').Value = '3.66'
$ws.Range('E42').Value = '  +1.11%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').Value = '130.51'
$ws.Range('E44').Value = '  +10.61%  '
$ws.Range('D45').Value = '0.0977'
$ws.Range('E45').Value = '  +1.85%  '
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = '0.0536'
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('D49').Value = '19.12'
$ws.Range('E49').Value = '  +1.28%  '
$ws.Range('E50').Value = '  +2.22%  '
$ws.Range('D51').Value = '1.943.64'
$ws.Range('E51').Value = '  -0.62%  '

# Restore the original style so the cells end up unstyled again, matching
# the source workbook (no explicit NumberFormat override left behind).
$dataRange.Style = "Normal"

